$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Fitness values (column C) from row 2 to row 252 are corrected to 7293
$ws.Range("C2:C252").Value = 7293
